# Auto-generated edit script: applies 2023-09-23 violent crime data update
# Updates column J (year 2023 running totals) and, for a few rows, column I (2022)
# across 47 worksheets, per the commit "Add data for 2023-09-23".

$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "Citywide Totals"; Cell = "J2"; OldValue = 5538; NewValue = 5565 }
    @{ Sheet = "Citywide Totals"; Cell = "J3"; OldValue = 5916; NewValue = 5936 }
    @{ Sheet = "Citywide Totals"; Cell = "I4"; OldValue = 1773; NewValue = 1775 }
    @{ Sheet = "Citywide Totals"; Cell = "J4"; OldValue = 1285; NewValue = 1291 }
    @{ Sheet = "Citywide Totals"; Cell = "J6"; OldValue = 7454; NewValue = 7504 }
    @{ Sheet = "Citywide Totals"; Cell = "I7"; OldValue = 26228; NewValue = 26230 }
    @{ Sheet = "Citywide Totals"; Cell = "J7"; OldValue = 20650; NewValue = 20753 }
    @{ Sheet = "Logan Square"; Cell = "J6"; OldValue = 180; NewValue = 182 }
    @{ Sheet = "Logan Square"; Cell = "J7"; OldValue = 283; NewValue = 285 }
    @{ Sheet = "Austin"; Cell = "J2"; OldValue = 357; NewValue = 358 }
    @{ Sheet = "Austin"; Cell = "J3"; OldValue = 395; NewValue = 396 }
    @{ Sheet = "Austin"; Cell = "J4"; OldValue = 74; NewValue = 75 }
    @{ Sheet = "Austin"; Cell = "J6"; OldValue = 436; NewValue = 437 }
    @{ Sheet = "Austin"; Cell = "J7"; OldValue = 1297; NewValue = 1301 }
    @{ Sheet = "Garfield Park"; Cell = "J3"; OldValue = 308; NewValue = 312 }
    @{ Sheet = "Garfield Park"; Cell = "J6"; OldValue = 325; NewValue = 328 }
    @{ Sheet = "Garfield Park"; Cell = "J7"; OldValue = 947; NewValue = 954 }
    @{ Sheet = "West Pullman"; Cell = "J2"; OldValue = 107; NewValue = 108 }
    @{ Sheet = "West Pullman"; Cell = "J7"; OldValue = 308; NewValue = 309 }
    @{ Sheet = "Grand Crossing"; Cell = "J2"; OldValue = 186; NewValue = 189 }
    @{ Sheet = "Grand Crossing"; Cell = "J6"; OldValue = 188; NewValue = 190 }
    @{ Sheet = "Grand Crossing"; Cell = "J7"; OldValue = 640; NewValue = 645 }
    @{ Sheet = "New City"; Cell = "J6"; OldValue = 186; NewValue = 185 }
    @{ Sheet = "New City"; Cell = "J7"; OldValue = 529; NewValue = 528 }
    @{ Sheet = "By Neighborhood"; Cell = "J2"; OldValue = 166; NewValue = 167 }
    @{ Sheet = "By Neighborhood"; Cell = "J4"; OldValue = 87; NewValue = 89 }
    @{ Sheet = "By Neighborhood"; Cell = "J5"; OldValue = 63; NewValue = 64 }
    @{ Sheet = "By Neighborhood"; Cell = "J7"; OldValue = 605; NewValue = 609 }
    @{ Sheet = "By Neighborhood"; Cell = "J8"; OldValue = 1297; NewValue = 1301 }
    @{ Sheet = "By Neighborhood"; Cell = "J9"; OldValue = 99; NewValue = 100 }
    @{ Sheet = "By Neighborhood"; Cell = "J10"; OldValue = 139; NewValue = 140 }
    @{ Sheet = "By Neighborhood"; Cell = "J11"; OldValue = 321; NewValue = 327 }
    @{ Sheet = "By Neighborhood"; Cell = "J12"; OldValue = 41; NewValue = 42 }
    @{ Sheet = "By Neighborhood"; Cell = "J14"; OldValue = 103; NewValue = 104 }
    @{ Sheet = "By Neighborhood"; Cell = "J16"; OldValue = 77; NewValue = 78 }
    @{ Sheet = "By Neighborhood"; Cell = "J18"; OldValue = 175; NewValue = 176 }
    @{ Sheet = "By Neighborhood"; Cell = "J19"; OldValue = 603; NewValue = 605 }
    @{ Sheet = "By Neighborhood"; Cell = "J20"; OldValue = 424; NewValue = 430 }
    @{ Sheet = "By Neighborhood"; Cell = "J23"; OldValue = 197; NewValue = 199 }
    @{ Sheet = "By Neighborhood"; Cell = "J25"; OldValue = 103; NewValue = 104 }
    @{ Sheet = "By Neighborhood"; Cell = "J27"; OldValue = 124; NewValue = 125 }
    @{ Sheet = "By Neighborhood"; Cell = "J29"; OldValue = 1156; NewValue = 1161 }
    @{ Sheet = "By Neighborhood"; Cell = "J31"; OldValue = 186; NewValue = 187 }
    @{ Sheet = "By Neighborhood"; Cell = "J33"; OldValue = 947; NewValue = 954 }
    @{ Sheet = "By Neighborhood"; Cell = "J36"; OldValue = 287; NewValue = 288 }
    @{ Sheet = "By Neighborhood"; Cell = "J37"; OldValue = 640; NewValue = 645 }
    @{ Sheet = "By Neighborhood"; Cell = "J40"; OldValue = 45; NewValue = 47 }
    @{ Sheet = "By Neighborhood"; Cell = "J41"; OldValue = 127; NewValue = 128 }
    @{ Sheet = "By Neighborhood"; Cell = "J42"; OldValue = 855; NewValue = 861 }
    @{ Sheet = "By Neighborhood"; Cell = "J43"; OldValue = 170; NewValue = 173 }
    @{ Sheet = "By Neighborhood"; Cell = "J51"; OldValue = 255; NewValue = 258 }
    @{ Sheet = "By Neighborhood"; Cell = "J52"; OldValue = 523; NewValue = 525 }
    @{ Sheet = "By Neighborhood"; Cell = "J53"; OldValue = 283; NewValue = 285 }
    @{ Sheet = "By Neighborhood"; Cell = "J54"; OldValue = 402; NewValue = 405 }
    @{ Sheet = "By Neighborhood"; Cell = "J55"; OldValue = 278; NewValue = 281 }
    @{ Sheet = "By Neighborhood"; Cell = "J60"; OldValue = 126; NewValue = 127 }
    @{ Sheet = "By Neighborhood"; Cell = "I63"; OldValue = 240; NewValue = 242 }
    @{ Sheet = "By Neighborhood"; Cell = "J63"; OldValue = 76; NewValue = 73 }
    @{ Sheet = "By Neighborhood"; Cell = "J65"; OldValue = 529; NewValue = 528 }
    @{ Sheet = "By Neighborhood"; Cell = "J67"; OldValue = 781; NewValue = 787 }
    @{ Sheet = "By Neighborhood"; Cell = "J73"; OldValue = 198; NewValue = 200 }
    @{ Sheet = "By Neighborhood"; Cell = "J76"; OldValue = 297; NewValue = 299 }
    @{ Sheet = "By Neighborhood"; Cell = "J78"; OldValue = 254; NewValue = 256 }
    @{ Sheet = "By Neighborhood"; Cell = "J79"; OldValue = 590; NewValue = 593 }
    @{ Sheet = "By Neighborhood"; Cell = "J84"; OldValue = 176; NewValue = 178 }
    @{ Sheet = "By Neighborhood"; Cell = "J85"; OldValue = 871; NewValue = 873 }
    @{ Sheet = "By Neighborhood"; Cell = "J88"; OldValue = 222; NewValue = 223 }
    @{ Sheet = "By Neighborhood"; Cell = "J90"; OldValue = 225; NewValue = 227 }
    @{ Sheet = "By Neighborhood"; Cell = "J91"; OldValue = 228; NewValue = 230 }
    @{ Sheet = "By Neighborhood"; Cell = "J95"; OldValue = 308; NewValue = 309 }
    @{ Sheet = "By Neighborhood"; Cell = "J96"; OldValue = 243; NewValue = 244 }
    @{ Sheet = "By Neighborhood"; Cell = "I101"; OldValue = 26228; NewValue = 26230 }
    @{ Sheet = "By Neighborhood"; Cell = "J101"; OldValue = 20650; NewValue = 20753 }
    @{ Sheet = "Gage Park"; Cell = "J3"; OldValue = 49; NewValue = 50 }
    @{ Sheet = "Gage Park"; Cell = "J7"; OldValue = 186; NewValue = 187 }
    @{ Sheet = "North Lawndale"; Cell = "J2"; OldValue = 193; NewValue = 195 }
    @{ Sheet = "North Lawndale"; Cell = "J6"; OldValue = 205; NewValue = 209 }
    @{ Sheet = "North Lawndale"; Cell = "J7"; OldValue = 781; NewValue = 787 }
    @{ Sheet = "South Deering"; Cell = "J2"; OldValue = 55; NewValue = 56 }
    @{ Sheet = "South Deering"; Cell = "J3"; OldValue = 57; NewValue = 58 }
    @{ Sheet = "South Deering"; Cell = "J7"; OldValue = 176; NewValue = 178 }
    @{ Sheet = "Loop"; Cell = "J2"; OldValue = 99; NewValue = 100 }
    @{ Sheet = "Loop"; Cell = "J6"; OldValue = 190; NewValue = 192 }
    @{ Sheet = "Loop"; Cell = "J7"; OldValue = 402; NewValue = 405 }
    @{ Sheet = "Englewood"; Cell = "J2"; OldValue = 345; NewValue = 347 }
    @{ Sheet = "Englewood"; Cell = "J3"; OldValue = 399; NewValue = 401 }
    @{ Sheet = "Englewood"; Cell = "J4"; OldValue = 65; NewValue = 66 }
    @{ Sheet = "Englewood"; Cell = "J7"; OldValue = 1156; NewValue = 1161 }
    @{ Sheet = "Chatham"; Cell = "J6"; OldValue = 222; NewValue = 224 }
    @{ Sheet = "Chatham"; Cell = "J7"; OldValue = 603; NewValue = 605 }
    @{ Sheet = "River North"; Cell = "J2"; OldValue = 44; NewValue = 45 }
    @{ Sheet = "River North"; Cell = "J6"; OldValue = 165; NewValue = 166 }
    @{ Sheet = "River North"; Cell = "J7"; OldValue = 297; NewValue = 299 }
    @{ Sheet = "Bridgeport"; Cell = "J6"; OldValue = 35; NewValue = 36 }
    @{ Sheet = "Bridgeport"; Cell = "J7"; OldValue = 103; NewValue = 104 }
    @{ Sheet = "Hermosa"; Cell = "J6"; OldValue = 68; NewValue = 69 }
    @{ Sheet = "Hermosa"; Cell = "J7"; OldValue = 127; NewValue = 128 }
    @{ Sheet = "Humboldt Park"; Cell = "J3"; OldValue = 171; NewValue = 173 }
    @{ Sheet = "Humboldt Park"; Cell = "J6"; OldValue = 438; NewValue = 442 }
    @{ Sheet = "Humboldt Park"; Cell = "J7"; OldValue = 855; NewValue = 861 }
    @{ Sheet = "Avondale"; Cell = "J6"; OldValue = 76; NewValue = 77 }
    @{ Sheet = "Avondale"; Cell = "J7"; OldValue = 139; NewValue = 140 }
    @{ Sheet = "Rogers Park"; Cell = "J2"; OldValue = 70; NewValue = 71 }
    @{ Sheet = "Rogers Park"; Cell = "J6"; OldValue = 70; NewValue = 71 }
    @{ Sheet = "Rogers Park"; Cell = "J7"; OldValue = 254; NewValue = 256 }
    @{ Sheet = "Lower West Side"; Cell = "J6"; OldValue = 137; NewValue = 140 }
    @{ Sheet = "Lower West Side"; Cell = "J7"; OldValue = 278; NewValue = 281 }
    @{ Sheet = "Douglas"; Cell = "J3"; OldValue = 67; NewValue = 68 }
    @{ Sheet = "Douglas"; Cell = "J6"; OldValue = 52; NewValue = 53 }
    @{ Sheet = "Douglas"; Cell = "J7"; OldValue = 197; NewValue = 199 }
    @{ Sheet = "West Ridge"; Cell = "J2"; OldValue = 72; NewValue = 73 }
    @{ Sheet = "West Ridge"; Cell = "J7"; OldValue = 243; NewValue = 244 }
    @{ Sheet = "Washington Park"; Cell = "J3"; OldValue = 93; NewValue = 94 }
    @{ Sheet = "Washington Park"; Cell = "J6"; OldValue = 51; NewValue = 52 }
    @{ Sheet = "Washington Park"; Cell = "J7"; OldValue = 228; NewValue = 230 }
    @{ Sheet = "Roseland"; Cell = "J2"; OldValue = 166; NewValue = 167 }
    @{ Sheet = "Roseland"; Cell = "J6"; OldValue = 165; NewValue = 167 }
    @{ Sheet = "Roseland"; Cell = "J7"; OldValue = 590; NewValue = 593 }
    @{ Sheet = "Chicago Lawn"; Cell = "J2"; OldValue = 117; NewValue = 120 }
    @{ Sheet = "Chicago Lawn"; Cell = "J3"; OldValue = 148; NewValue = 150 }
    @{ Sheet = "Chicago Lawn"; Cell = "J4"; OldValue = 38; NewValue = 39 }
    @{ Sheet = "Chicago Lawn"; Cell = "J7"; OldValue = 424; NewValue = 430 }
    @{ Sheet = "Calumet Heights"; Cell = "J2"; OldValue = 47; NewValue = 48 }
    @{ Sheet = "Calumet Heights"; Cell = "J7"; OldValue = 175; NewValue = 176 }
    @{ Sheet = "Grand Boulevard"; Cell = "J3"; OldValue = 92; NewValue = 93 }
    @{ Sheet = "Grand Boulevard"; Cell = "J7"; OldValue = 287; NewValue = 288 }
    @{ Sheet = "Auburn Gresham"; Cell = "J3"; OldValue = 182; NewValue = 183 }
    @{ Sheet = "Auburn Gresham"; Cell = "J4"; OldValue = 24; NewValue = 26 }
    @{ Sheet = "Auburn Gresham"; Cell = "J6"; OldValue = 196; NewValue = 197 }
    @{ Sheet = "Auburn Gresham"; Cell = "J7"; OldValue = 605; NewValue = 609 }
    @{ Sheet = "East Side"; Cell = "J2"; OldValue = 46; NewValue = 47 }
    @{ Sheet = "East Side"; Cell = "J7"; OldValue = 103; NewValue = 104 }
    @{ Sheet = "Belmont Cragin"; Cell = "J2"; OldValue = 98; NewValue = 99 }
    @{ Sheet = "Belmont Cragin"; Cell = "J5"; OldValue = 5; NewValue = 6 }
    @{ Sheet = "Belmont Cragin"; Cell = "J6"; OldValue = 131; NewValue = 135 }
    @{ Sheet = "Belmont Cragin"; Cell = "J7"; OldValue = 321; NewValue = 327 }
    @{ Sheet = "Avalon Park"; Cell = "J2"; OldValue = 26; NewValue = 27 }
    @{ Sheet = "Avalon Park"; Cell = "J7"; OldValue = 99; NewValue = 100 }
    @{ Sheet = "Portage Park"; Cell = "J2"; OldValue = 68; NewValue = 69 }
    @{ Sheet = "Portage Park"; Cell = "J6"; OldValue = 65; NewValue = 66 }
    @{ Sheet = "Portage Park"; Cell = "J7"; OldValue = 198; NewValue = 200 }
    @{ Sheet = "Albany Park"; Cell = "J6"; OldValue = 62; NewValue = 63 }
    @{ Sheet = "Albany Park"; Cell = "J7"; OldValue = 166; NewValue = 167 }
    @{ Sheet = "United Center"; Cell = "J6"; OldValue = 101; NewValue = 102 }
    @{ Sheet = "United Center"; Cell = "J7"; OldValue = 222; NewValue = 223 }
    @{ Sheet = "Armour Square"; Cell = "J2"; OldValue = 21; NewValue = 22 }
    @{ Sheet = "Armour Square"; Cell = "J7"; OldValue = 63; NewValue = 64 }
    @{ Sheet = "Edgewater"; Cell = "J6"; OldValue = 42; NewValue = 43 }
    @{ Sheet = "Edgewater"; Cell = "J7"; OldValue = 124; NewValue = 125 }
    @{ Sheet = "Washington Heights"; Cell = "J3"; OldValue = 65; NewValue = 66 }
    @{ Sheet = "Washington Heights"; Cell = "J6"; OldValue = 64; NewValue = 65 }
    @{ Sheet = "Washington Heights"; Cell = "J7"; OldValue = 225; NewValue = 227 }
    @{ Sheet = "Little Italy, UIC"; Cell = "J6"; OldValue = 96; NewValue = 99 }
    @{ Sheet = "Little Italy, UIC"; Cell = "J7"; OldValue = 255; NewValue = 258 }
    @{ Sheet = "Morgan Park"; Cell = "J6"; OldValue = 36; NewValue = 37 }
    @{ Sheet = "Morgan Park"; Cell = "J7"; OldValue = 126; NewValue = 127 }
    @{ Sheet = "Hyde Park"; Cell = "J3"; OldValue = 34; NewValue = 36 }
    @{ Sheet = "Hyde Park"; Cell = "J6"; OldValue = 98; NewValue = 99 }
    @{ Sheet = "Hyde Park"; Cell = "J7"; OldValue = 170; NewValue = 173 }
    @{ Sheet = "South Shore"; Cell = "J3"; OldValue = 317; NewValue = 318 }
    @{ Sheet = "South Shore"; Cell = "J6"; OldValue = 252; NewValue = 253 }
    @{ Sheet = "South Shore"; Cell = "J7"; OldValue = 871; NewValue = 873 }
    @{ Sheet = "Hegewisch"; Cell = "J2"; OldValue = 18; NewValue = 20 }
    @{ Sheet = "Hegewisch"; Cell = "J7"; OldValue = 45; NewValue = 47 }
    @{ Sheet = "Little Village"; Cell = "J6"; OldValue = 209; NewValue = 211 }
    @{ Sheet = "Little Village"; Cell = "J7"; OldValue = 523; NewValue = 525 }
    @{ Sheet = "Archer Heights"; Cell = "J6"; OldValue = 31; NewValue = 33 }
    @{ Sheet = "Archer Heights"; Cell = "J7"; OldValue = 87; NewValue = 89 }
    @{ Sheet = "Beverly"; Cell = "J6"; OldValue = 27; NewValue = 28 }
    @{ Sheet = "Beverly"; Cell = "J7"; OldValue = 41; NewValue = 42 }
    @{ Sheet = "Bucktown"; Cell = "J6"; OldValue = 60; NewValue = 61 }
    @{ Sheet = "Bucktown"; Cell = "J7"; OldValue = 77; NewValue = 78 }
)

foreach ($chg in $changes) {
    $ws = $wb.Worksheets.Item($chg.Sheet)
    $ws.Range($chg.Cell).Value = $chg.NewValue
}

Write-Output ("Applied {0} cell updates across {1} worksheets." -f $changes.Count, ($changes | Select-Object -ExpandProperty Sheet -Unique).Count)
